# Add 2022-Q1 data
# 1) Insert a new "2022-Q1" sheet (cloned from "2021-Q4" so header text/styles match)
#    positioned right before "总计".
# 2) Replace its per-fund rows with the new 2022-Q1 holdings.
# 3) Prepend a "2022-Q1" summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet by cloning "2021-Q4" (same header/styles)
# and placing it immediately before "总计".
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template had 59 data rows (rows 2-60); the new sheet only needs 29
# (rows 2-30), so drop the extra rows entirely (shrinks the used range too).
$newSheet.Rows("31:60").Delete()

# ---------------------------------------------------------------------------
# Step 2: fund-level holdings for 2022-Q1
# Columns: A idx(auto) | B code | C name | D scale | E total position |
#          F position ratio | G market value (亿元) | H position rank
# ---------------------------------------------------------------------------
$fundData = @(
    @('003961', '易方达瑞程灵活配置混合A', '34.99', '91.05', '2.96', '1.0357', 8),
    @('012647', '中欧洞见一年持有混合', '33.02', '65.46', '2.64', '0.8717', 8),
    @('010330', '东吴兴享成长混合A', '11.63', '80.15', '4.63', '0.5385', 6),
    @('003962', '易方达瑞程灵活配置混合C', '9.83', '91.05', '2.96', '0.2910', 8),
    @('013004', '国泰价值领航股票A', '7.51', '94.17', '3.43', '0.2576', 9),
    @('006449', '浙商汇金量化精选灵活配置混合', '4.27', '74.69', '3.89', '0.1661', 2),
    @('200010', '长城双动力混合', '2.03', '86.77', '3.70', '0.0751', 5),
    @('006385', '华泰保兴研究智选灵活配置混合A', '2.15', '83.65', '3.11', '0.0669', 8),
    @('580006', '东吴新经济混合', '1.17', '91.56', '4.78', '0.0559', 7),
    @('005933', '新疆前海联合先进制造灵活配置混合A', '1.26', '89.79', '4.29', '0.0541', 10),
    @('005571', '中银证券新能源灵活配置混合A', '0.91', '90.25', '5.08', '0.0462', 6),
    @('012696', '同泰数字经济主题股票型证券投资基金A', '1.54', '93.98', '2.95', '0.0454', 10),
    @('013204', '恒生前海恒源天利债A', '1.29', '21.66', '2.11', '0.0272', 2),
    @('011462', '东吴兴享成长混合C', '0.33', '80.15', '4.63', '0.0153', 6),
    @('007439', '东海科技动力混合A', '0.36', '90.01', '4.20', '0.0151', 8),
    @('005572', '中银证券新能源灵活配置混合C', '0.28', '90.25', '5.08', '0.0142', 6),
    @('012697', '同泰数字经济主题股票型证券投资基金C', '0.44', '93.98', '2.95', '0.0130', 10),
    @('005169', '华泰保兴策略精选灵活配置混合A', '0.38', '84.18', '3.17', '0.0120', 10),
    @('013005', '国泰价值领航股票C', '0.33', '94.17', '3.43', '0.0113', 9),
    @('007463', '东海科技动力混合C', '0.20', '90.01', '4.20', '0.0084', 8),
    @('005170', '华泰保兴策略精选灵活配置混合C', '0.23', '84.18', '3.17', '0.0073', 10),
    @('005934', '新疆前海联合先进制造灵活配置混合C', '0.09', '89.79', '4.29', '0.0039', 10),
    @('001351', '诺安中证500指数增强A', '0.48', '94.42', '0.70', '0.0034', 1),
    @('006538', '东海核心价值精选混合', '0.03', '89.24', '4.22', '0.0013', 8),
    @('006386', '华泰保兴研究智选灵活配置混合C', '0.02', '83.65', '3.11', '0.0006', 8),
    @('010355', '诺安中证500指数增强C', '0.04', '94.42', '0.70', '0.0003', 1),
    @('004005', '东方民丰回报赢安混合A', '0.01', '28.31', '1.68', '0.0002', 9),
    @('004006', '东方民丰回报赢安混合C', '0.00', '28.31', '1.68', '0', 9),
    @('013205', '恒生前海恒源天利债C', '0.00', '21.66', '2.11', '0', 2)
)

# Text-like columns must keep leading zeros / trailing zeros verbatim
# (e.g. fund code "003961", scale "34.99"), so force them to text before
# assignment -- otherwise Excel would auto-coerce them to numbers.
$newSheet.Range("B2:G30").NumberFormat = "@"

$r = 2
foreach ($item in $fundData) {
    $newSheet.Cells.Item($r, 2).Value = $item[0]
    $newSheet.Cells.Item($r, 3).Value = $item[1]
    $newSheet.Cells.Item($r, 4).Value = $item[2]
    $newSheet.Cells.Item($r, 5).Value = $item[3]
    $newSheet.Cells.Item($r, 6).Value = $item[4]

    if ($item[5] -eq '0') {
        # Zero market value is stored as a genuine number, not text.
        $newSheet.Cells.Item($r, 7).NumberFormat = "General"
        $newSheet.Cells.Item($r, 7).Value = 0
    } else {
        $newSheet.Cells.Item($r, 7).Value = $item[5]
    }

    $newSheet.Cells.Item($r, 8).Value = $item[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Step 3: prepend the 2022-Q1 summary row to "总计"
# ---------------------------------------------------------------------------
$totalRows = @(
    @('2022-Q1', 29, 3.64),
    @('2021-Q4', 59, 6.33),
    @('2021-Q3', 36, 3.32),
    @('2021-Q2', 20, 1.9),
    @('2021-Q1', 12, 1.28),
    @('2020-Q4', 8, 0.84)
)

# Row 7 is brand new (the table only had 6 rows before); give its A cell the
# same style (bold + boxed) the other index cells already carry.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

$totalSheet.Range("A2:D7").ClearContents()

$r = 2
$idx = 0
foreach ($item in $totalRows) {
    $totalSheet.Cells.Item($r, 1).Value = $idx
    $totalSheet.Cells.Item($r, 2).Value = $item[0]
    $totalSheet.Cells.Item($r, 3).Value = $item[1]
    $totalSheet.Cells.Item($r, 4).Value = $item[2]
    $r = $r + 1
    $idx = $idx + 1
}

Write-Host "2022-Q1 sheet added and 总计 updated."
